$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B2").Value = 660
$ws.Range("C2").Value = 1946
$ws.Range("D2").Value = 4941

$ws.Range("B5").Value = 184
$ws.Range("C5").Value = 181
$ws.Range("D5").Value = 485

$ws.Range("B13").Value = 5792
$ws.Range("C13").Value = 16201
$ws.Range("D13").Value = 44863

$ws.Range("B15").Value = 79
$ws.Range("C15").Value = 96
$ws.Range("D15").Value = 113

$ws.Range("B16").Select()
$excel.ActiveWindow.ScrollRow = 2
